$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 'profile-thread-7388a4a3-91b5-42a4-a1e2-7b295561673c'
$ws.Range("A3").Value = 'https://www.uidaho.edu/cals/people/mark-mcguire'
$ws.Range("C3").Value = '{''exception_type'': ''HTTPError'', ''message'': ''404 Client Error: Not Found for url: https://www.uidaho.edu/cals/people/mark-mcguire'', ''status_code'': 404, ''traceback'': ''Traceback (most recent call last):\n  File "/Users/colesummers/Documents/GitHub/profile-extractor/src/nodes.py", line 88, in fetch_html\n    response.raise_for_status()  # Raise HTTPError for bad responses (4xx or 5xx)\n  File "/Users/colesummers/Documents/GitHub/profile-extractor/.venv/lib/python3.9/site-packages/requests/models.py", line 1024, in raise_for_status\n    raise HTTPError(http_error_msg, response=self)\nrequests.exceptions.HTTPError: 404 Client Error: Not Found for url: https://www.uidaho.edu/cals/people/mark-mcguire\n''}'
$ws.Range("D3").Value = 'profile-thread-149f4e6a-3b71-42e1-8e2d-781b89405e91'
$ws.Range("A4").Value = 'https://www.uidaho.edu/cals/agricultural-economics-and-rural-sociology/our-people/paul-lewin'
$ws.Range("B4").Value = 'Failed to fetch URL'
$ws.Range("C4").Value = '{''exception_type'': ''HTTPError'', ''message'': ''404 Client Error: Not Found for url: https://www.uidaho.edu/cals/agricultural-economics-and-rural-sociology/our-people/paul-lewin'', ''status_code'': 404, ''traceback'': ''Traceback (most recent call last):\n  File "/Users/colesummers/Documents/GitHub/profile-extractor/src/nodes.py", line 88, in fetch_html\n    response.raise_for_status()  # Raise HTTPError for bad responses (4xx or 5xx)\n  File "/Users/colesummers/Documents/GitHub/profile-extractor/.venv/lib/python3.9/site-packages/requests/models.py", line 1024, in raise_for_status\n    raise HTTPError(http_error_msg, response=self)\nrequests.exceptions.HTTPError: 404 Client Error: Not Found for url: https://www.uidaho.edu/cals/agricultural-economics-and-rural-sociology/our-people/paul-lewin\n''}'
$ws.Range("D4").Value = 'profile-thread-b3a8e2da-4b0c-49bd-a719-044545da9d42'
$ws.Range("A5").Value = 'https://www.uidaho.edu/cals/animal-veterinary-and-food-sciences/our-people/andrzej-paszczynski'
$ws.Range("B5").Value = 'Failed to fetch URL'
$ws.Range("C5").Value = '{''exception_type'': ''HTTPError'', ''message'': ''404 Client Error: Not Found for url: https://www.uidaho.edu/cals/animal-veterinary-and-food-sciences/our-people/andrzej-paszczynski'', ''status_code'': 404, ''traceback'': ''Traceback (most recent call last):\n  File "/Users/colesummers/Documents/GitHub/profile-extractor/src/nodes.py", line 88, in fetch_html\n    response.raise_for_status()  # Raise HTTPError for bad responses (4xx or 5xx)\n  File "/Users/colesummers/Documents/GitHub/profile-extractor/.venv/lib/python3.9/site-packages/requests/models.py", line 1024, in raise_for_status\n    raise HTTPError(http_error_msg, response=self)\nrequests.exceptions.HTTPError: 404 Client Error: Not Found for url: https://www.uidaho.edu/cals/animal-veterinary-and-food-sciences/our-people/andrzej-paszczynski\n''}'
$ws.Range("D5").Value = 'profile-thread-6d6897ba-da41-42a9-b5bf-0b00e4783c79'
$ws.Range("A6").Value = 'https://www.uidaho.edu/cals/animal-veterinary-and-food-sciences/our-people/barbara-nielsen'
$ws.Range("B6").Value = 'Failed to fetch URL'
$ws.Range("C6").Value = '{''exception_type'': ''HTTPError'', ''message'': ''404 Client Error: Not Found for url: https://www.uidaho.edu/cals/animal-veterinary-and-food-sciences/our-people/barbara-nielsen'', ''status_code'': 404, ''traceback'': ''Traceback (most recent call last):\n  File "/Users/colesummers/Documents/GitHub/profile-extractor/src/nodes.py", line 88, in fetch_html\n    response.raise_for_status()  # Raise HTTPError for bad responses (4xx or 5xx)\n  File "/Users/colesummers/Documents/GitHub/profile-extractor/.venv/lib/python3.9/site-packages/requests/models.py", line 1024, in raise_for_status\n    raise HTTPError(http_error_msg, response=self)\nrequests.exceptions.HTTPError: 404 Client Error: Not Found for url: https://www.uidaho.edu/cals/animal-veterinary-and-food-sciences/our-people/barbara-nielsen\n''}'
$ws.Range("D6").Value = 'profile-thread-0049fa1d-9edd-4f45-ad28-aef5bef13f4f'
